$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 1150
$ws.Range("J33").Value = 39
$ws.Range("L33").Value = 39
$ws.Range("N33").Value = -497
# row 40
$ws.Range("H40").Value = 7605.5884
$ws.Range("I40").Value = 6500.75
$ws.Range("J40").Value = 7945.5386
$ws.Range("K40").Value = 6500.75
$ws.Range("L40").Value = 7945.5386
$ws.Range("M40").Value = -6325.75
$ws.Range("N40").Value = -8295.5386
# row 87
$ws.Range("H87").Value = 89500
$ws.Range("I87").Value = 59000
$ws.Range("J87").Value = 92272.73
$ws.Range("K87").Value = 59000
$ws.Range("L87").Value = 92272.73
$ws.Range("M87").Value = -57752
$ws.Range("N87").Value = -94768.73
# row 90
$ws.Range("H90").Value = 89500
$ws.Range("I90").Value = 59000
$ws.Range("J90").Value = 92272.73
$ws.Range("K90").Value = 177000
$ws.Range("L90").Value = 276818.19
$ws.Range("M90").Value = -170760
$ws.Range("N90").Value = -289298.19
# row 116
$ws.Range("H116").Value = 7760.4
$ws.Range("I116").Value = 3479.3157
$ws.Range("K116").Value = 3479.3157
$ws.Range("M116").Value = -37.31570000000011
# row 137
$ws.Range("H137").Value = 647961.6
$ws.Range("I137").Value = 478212.44
$ws.Range("J137").Value = 1004434.9
$ws.Range("K137").Value = 1434637.32
$ws.Range("L137").Value = 3013304.7
$ws.Range("M137").Value = -1432087.32
$ws.Range("N137").Value = -3018404.7

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 3
$ws.Range("H3").Value = 55000
$ws.Range("I3").Value = 55000
$ws.Range("K3").Value = 55000
$ws.Range("M3").Value = -54885
# row 32
$ws.Range("H32").Value = 4754.98
$ws.Range("I32").Value = 4754.98
$ws.Range("K32").Value = 4754.98
$ws.Range("M32").Value = -4467.98
# row 74
$ws.Range("H74").Value = 2073.476
$ws.Range("I74").Value = 1752.4445
$ws.Range("J74").Value = 3999.6667
$ws.Range("K74").Value = 1752.4445
$ws.Range("L74").Value = 3999.6667
$ws.Range("M74").Value = -878.4445000000001
$ws.Range("N74").Value = -5747.6667
# row 77
$ws.Range("H77").Value = 2073.476
$ws.Range("I77").Value = 1752.4445
$ws.Range("J77").Value = 3999.6667
$ws.Range("K77").Value = 8762.2225
$ws.Range("L77").Value = 19998.3335
$ws.Range("M77").Value = -4394.2225
$ws.Range("N77").Value = -28734.3335

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 14
$ws.Range("H14").Value = 1166.6666
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 1500
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = -828
$ws.Range("N14").Value = -1844
# row 17
$ws.Range("H17").Value = 2599.8
$ws.Range("J17").Value = 2599.8
$ws.Range("L17").Value = 2599.8
$ws.Range("N17").Value = -2943.8
# row 22
$ws.Range("H22").Value = 297.85715
$ws.Range("I22").Value = 277
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 277
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -104
$ws.Range("N22").Value = -696
# row 132
$ws.Range("H132").Value = 50000
$ws.Range("J132").Value = 50000
$ws.Range("L132").Value = 50000
$ws.Range("N132").Value = -60120

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1079.25
$ws.Range("I16").Value = 872
$ws.Range("K16").Value = 872
$ws.Range("M16").Value = -585
# row 68
$ws.Range("H68").Value = 100000
$ws.Range("J68").Value = 120000
$ws.Range("L68").Value = 120000
$ws.Range("N68").Value = -121498
# row 71
$ws.Range("H71").Value = 100000
$ws.Range("J71").Value = 120000
$ws.Range("L71").Value = 360000
$ws.Range("N71").Value = -367488
# row 99
$ws.Range("H99").Value = 5305.727
$ws.Range("I99").Value = 3795.6365
$ws.Range("J99").Value = 6815.8184
$ws.Range("K99").Value = 3795.6365
$ws.Range("L99").Value = 6815.8184
$ws.Range("M99").Value = -2297.6365
$ws.Range("N99").Value = -9811.8184
# row 113
$ws.Range("H113").Value = 1079.25
$ws.Range("I113").Value = 872
$ws.Range("K113").Value = 872
$ws.Range("M113").Value = 1298
# row 126
$ws.Range("H126").Value = 5305.727
$ws.Range("I126").Value = 3795.6365
$ws.Range("J126").Value = 6815.8184
$ws.Range("K126").Value = 11386.9095
$ws.Range("L126").Value = 20447.4552
$ws.Range("M126").Value = -8916.9095
$ws.Range("N126").Value = -25387.4552

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 68.07692
$ws.Range("I2").Value = 46.8
$ws.Range("K2").Value = 280.8
$ws.Range("M2").Value = -167.8
# row 5
$ws.Range("H5").Value = 31623.076
$ws.Range("I5").Value = 50287.812
$ws.Range("J5").Value = 1759.5
$ws.Range("K5").Value = 150863.436
$ws.Range("L5").Value = 5278.5
$ws.Range("M5").Value = -150751.436
$ws.Range("N5").Value = -5502.5
# row 42
$ws.Range("H42").Value = 5000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = ""
# row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""
# row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""
# row 98
$ws.Range("H98").Value = 3006
$ws.Range("J98").Value = 2031.3334
$ws.Range("L98").Value = 6094.0002
$ws.Range("N98").Value = -9090.0002
# row 122
$ws.Range("H122").Value = 50840.4
$ws.Range("I122").Value = 723.75
$ws.Range("K122").Value = 6513.75
$ws.Range("M122").Value = -4063.75
# row 132
$ws.Range("H132").Value = 5938.25
$ws.Range("I132").Value = 4269.3335
$ws.Range("J132").Value = 6939.6
$ws.Range("K132").Value = 38424.0015
$ws.Range("L132").Value = 62456.4
$ws.Range("M132").Value = -35894.0015
$ws.Range("N132").Value = -67516.39999999999
# row 135
$ws.Range("H135").Value = 31623.076
$ws.Range("I135").Value = 50287.812
$ws.Range("J135").Value = 1759.5
$ws.Range("K135").Value = 452590.308
$ws.Range("L135").Value = 15835.5
$ws.Range("M135").Value = -450055.308
$ws.Range("N135").Value = -20905.5
# row 137
$ws.Range("H137").Value = 2637.739
$ws.Range("I137").Value = 1589
$ws.Range("K137").Value = 4767
$ws.Range("M137").Value = 333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 2159.1428
$ws.Range("I102").Value = 499.66666
$ws.Range("J102").Value = 2435.7222
$ws.Range("K102").Value = 499.66666
$ws.Range("L102").Value = 2435.7222
$ws.Range("M102").Value = 1122.33334
$ws.Range("N102").Value = -5679.7222
# row 113
$ws.Range("H113").Value = 10112.75
$ws.Range("I113").Value = 3509.875
$ws.Range("K113").Value = 3509.875
$ws.Range("M113").Value = -1339.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 349645.53
$ws.Range("I7").Value = 4962.6875
$ws.Range("J7").Value = 773870.5600000001
$ws.Range("K7").Value = 4962.6875
$ws.Range("L7").Value = 773870.5600000001
$ws.Range("M7").Value = -4850.6875
$ws.Range("N7").Value = -774094.5600000001
# row 40
$ws.Range("H40").Value = 2176930.5
$ws.Range("I40").Value = 2780279.2
$ws.Range("K40").Value = 2780279.2
$ws.Range("M40").Value = -2780143.2
# row 46
$ws.Range("H46").Value = 3069.0435
$ws.Range("I46").Value = 2530.077
$ws.Range("J46").Value = 3769.7
$ws.Range("K46").Value = 2530.077
$ws.Range("L46").Value = 3769.7
$ws.Range("M46").Value = -2342.077
$ws.Range("N46").Value = -4145.7
# row 122
$ws.Range("H122").Value = 456551.38
$ws.Range("I122").Value = 2019.1765
$ws.Range("J122").Value = 2001960.8
$ws.Range("K122").Value = 6057.529500000001
$ws.Range("L122").Value = 6005882.4
$ws.Range("M122").Value = -3607.529500000001
$ws.Range("N122").Value = -6010782.4
# row 126
$ws.Range("H126").Value = 349645.53
$ws.Range("I126").Value = 4962.6875
$ws.Range("J126").Value = 773870.5600000001
$ws.Range("K126").Value = 14888.0625
$ws.Range("L126").Value = 2321611.68
$ws.Range("M126").Value = -12418.0625
$ws.Range("N126").Value = -2326551.68

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 41
$ws.Range("H41").Value = 7710.222
$ws.Range("I41").Value = 2498
$ws.Range("J41").Value = 8361.75
$ws.Range("K41").Value = 2498
$ws.Range("L41").Value = 8361.75
$ws.Range("M41").Value = -2108
$ws.Range("N41").Value = -9141.75
# row 100
$ws.Range("H100").Value = 1007
$ws.Range("I100").Value = 1112.7858
$ws.Range("K100").Value = 2225.5716
$ws.Range("M100").Value = -1684.5716
# row 132
$ws.Range("H132").Value = 17824.096
$ws.Range("I132").Value = 1440.9454
$ws.Range("J132").Value = 130458.25
$ws.Range("K132").Value = 4322.8362
$ws.Range("L132").Value = 391374.75
$ws.Range("M132").Value = -1792.8362
$ws.Range("N132").Value = -396434.75
